$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Update sheet (tab) name to reflect new "through" date
$ws.Name = "Through 2022-03-01"

# Update the header label in I1 (shared string "2022 (through 02-28)" -> "2022 (through 03-01)")
$ws.Range("I1").Value = "2022 (through 03-01)"

# Add new March 2022 data point
$ws.Range("I4").Value = 8

# Update the yearly total to include the new March value (300 -> 308)
$ws.Range("I14").Value = 308
